# Auto commit at 2025-10-16  7:48:37.02
#
# Updates the "Metrics" sheet's raw metric values (B2:B13), which flow
# through formulas into the "today" sheet (B11:B22, E11:E22, F11:F22),
# and the day-rollover formula on "today"!A1 (=TODAY()-1, volatile).
# Also moves the active sheet/selection from Metrics!H22 to today!H16,
# leaving Metrics' own remembered selection at C21.

$wb = $excel.ActiveWorkbook

# --- Update the raw metrics on the "Metrics" sheet -----------------------
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 208105.08
$metrics.Range("B3").Value  = 171070.46000000002
$metrics.Range("B4").Value  = 66154.720000000001
$metrics.Range("B5").Value  = 8240
$metrics.Range("B6").Value  = 4575236.55
$metrics.Range("B7").Value  = 3860889.13
$metrics.Range("B8").Value  = 1336756.8600000003
$metrics.Range("B9").Value  = 177241
$metrics.Range("B10").Value = 33040560.350999828
$metrics.Range("B11").Value = 19890759.200000003
$metrics.Range("B12").Value = 11618465.750000002
$metrics.Range("B13").Value = 1274868

# Leave the Metrics sheet's remembered selection at C21 (its cursor
# position before the user switched away to the "today" sheet).
$metrics.Range("C21").Select()

# --- Recalculate so dependent formulas (today!, TODAY()-1, etc.) refresh -
$excel.Calculate()

# --- Switch the active sheet/selection to "today" ------------------------
$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("H16").Select()
